$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 70 ---
$ws.Cells.Item(70, 1).Value = 44704
$ws.Cells.Item(70, 2).Value = 0.40972222222222227
$ws.Cells.Item(70, 3).Value = 0.48958333333333331
$ws.Cells.Item(70, 5).Value = "Documentation"
$ws.Cells.Item(70, 6).Value = "Révision de la structure et orthographe, mise en page et partie manquante/remise à plus tard"

# --- Row 71 ---
$ws.Cells.Item(71, 1).Value = 44704
$ws.Cells.Item(71, 2).Value = 0.48958333333333331
$ws.Cells.Item(71, 3).Value = 0.51041666666666663
$ws.Cells.Item(71, 5).Value = "Test par des tier"
$ws.Cells.Item(71, 6).Value = "Mise en place de l'environnement + Build`nTesteurs:`nAntoine Dubois`nGaetan Epars"
$ws.Cells.Item(71, 7).Value = "Quelques tests sont mals formulés d'autres sont carrément faut dans le résultats attendu"

# --- Row 72 ---
$ws.Cells.Item(72, 1).Value = 44704
$ws.Cells.Item(72, 2).Value = 0.5625
$ws.Cells.Item(72, 3).Value = 0.57638888888888895
$ws.Cells.Item(72, 5).Value = "Dossier de réalisation"

# --- Row 73 ---
$ws.Cells.Item(73, 1).Value = 44704
$ws.Cells.Item(73, 2).Value = 0.57638888888888895
$ws.Cells.Item(73, 3).Value = 0.59375
$ws.Cells.Item(73, 5).Value = "Entretien avec chef de projet"

# --- Row 74 ---
$ws.Cells.Item(74, 1).Value = 44704
$ws.Cells.Item(74, 2).Value = 0.59375
$ws.Cells.Item(74, 3).Value = 0.62847222222222221
$ws.Cells.Item(74, 5).Value = "Dossier de réalisation"

# --- Finish Description/Solutions cells in the order that matches the original authoring ---
$ws.Cells.Item(73, 6).Value = "Il faut lier le pathfinding avec le system de mouvement actuelle! Et faire un système de tir!"
$ws.Cells.Item(72, 6).Value = "Partie régiment"

# --- Update selection / view to match final state ---
$ws.Range("F74").Select() | Out-Null
